# The chart on the single slide has its chart-level title (the "Geomean /
# BOP / DA-AMPM / SPP / PPF" text block) removed, turning off the chart
# title entirely (c:title element dropped, autoTitleDeleted flipped on).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$chartShape = $s.Shapes.Item(1)
$chart = $chartShape.Chart

$chart.HasTitle = $false
